{"js": "// The document had a module heading \"\u67e5\u8be2\u4f59\u989d\" (\"Query balance\") that the\n// author edited in place to \"\u786e\u8ba4\u4f59\u989d\" (\"Confirm balance\") by selecting the\n// first two characters (\"\u67e5\u8be2\") and typing \"\u786e\u8ba4\" over them. Word leaves the\n// trailing \"\u4f59\u989d\" run untouched and relocates the internal \"_GoBack\"\n// bookmark (marks the last edit point) from the end of the document to the\n// spot right after the newly typed text.\n\nconst body = context.document.body;\n\n// 1) The \"_GoBack\" bookmark currently sits at the very end of the document\n//    (after the last embedded drawing). Remove it from there first.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the unique heading text \"\u67e5\u8be2\u4f59\u989d\" and narrow down to just the\n//    \"\u67e5\u8be2\" portion that gets replaced.\nconst headingRange = body.search(\"\u67e5\u8be2\u4f59\u989d\", { matchCase: true }).getFirstOrNullObject();\nheadingRange.load(\"text\");\nawait context.sync();\n\nif (!headingRange.isNullObject) {\n  const targetRange = headingRange\n    .search(\"\u67e5\u8be2\", { matchCase: true })\n    .getFirstOrNullObject();\n  targetRange.load(\"text\");\n  await context.sync();\n\n  if (!targetRange.isNullObject) {\n    // 3) Replace \"\u67e5\u8be2\" with \"\u786e\u8ba4\" in place.\n    targetRange.insertText(\"\u786e\u8ba4\", Word.InsertLocation.replace);\n    await context.sync();\n\n    // 4) Re-insert the \"_GoBack\" bookmark right after the replaced text,\n    //    before the remaining \"\u4f59\u989d\" text, matching Word's own behaviour\n    //    when tracking the last edited location.\n    const insertionPoint = targetRange.getRange(Word.RangeLocation.end);\n    insertionPoint.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# The document had a module heading \"\u67e5\u8be2\u4f59\u989d\" (\"Query balance\") that the\n# author edited in place to \"\u786e\u8ba4\u4f59\u989d\" (\"Confirm balance\") by selecting the\n# first two characters (\"\u67e5\u8be2\") and typing \"\u786e\u8ba4\" over them. Word leaves the\n# trailing \"\u4f59\u989d\" text untouched and relocates the internal \"_GoBack\"\n# bookmark (marks the last edit point) from the end of the document to the\n# spot right after the newly typed text.\n\n$d = $word.ActiveDocument\n\n# 1) The \"_GoBack\" bookmark currently sits at the very end of the document\n#    (after the last embedded drawing). Remove it from there first.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate the unique heading text \"\u67e5\u8be2\u4f59\u989d\".\n$full = $d.Content\n$found = $full.Find.Execute(\"\u67e5\u8be2\u4f59\u989d\")\n\nif ($found) {\n    $start = $full.Start\n\n    # 3) Drop a throwaway bookmark right before \"\u67e5\u8be2\" so the subsequent\n    #    text replace can't silently merge back into the preceding \" \" run\n    #    (mirrors the run boundary Word leaves after a real in-place edit).\n    $leftRange = $d.Range($start, $start)\n    $d.Bookmarks.Add(\"_TempSplit\", $leftRange)\n\n    # 4) Bookmark the gap between \"\u67e5\u8be2\" and \"\u4f59\u989d\" -- this becomes the\n    #    relocated \"_GoBack\" bookmark.\n    $bmRange = $d.Range($start + 2, $start + 2)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n    # 5) Replace \"\u67e5\u8be2\" with \"\u786e\u8ba4\" in place.\n    $target = $d.Range($start, $start + 2)\n    $target.Text = \"\u786e\u8ba4\"\n\n    # 6) Remove the throwaway helper bookmark.\n    $d.Bookmarks(\"_TempSplit\").Delete()\n}\n"}
